$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix rows whose home/away matchup + odds data had been shuffled between
#        neighbouring fixtures played on the same matchday (index/date columns A-E stay put) ---

# Row 15
$ws.Range('F15').Value = 'Orijent'
$ws.Range('G15').Value = 0
$ws.Range('H15').Value = 'Vukovar 1991'
$ws.Range('I15').Value = 1
$ws.Range('J15').Value = 2.76
$ws.Range('K15').Value = '25/08/2023 05:12'
$ws.Range('L15').Value = 1.93
$ws.Range('M15').Value = '26/08/2023 16:55'
$ws.Range('N15').Value = 3.12
$ws.Range('O15').Value = '25/08/2023 05:12'
$ws.Range('P15').Value = 3.54
$ws.Range('Q15').Value = '26/08/2023 16:52'
$ws.Range('R15').Value = 2.31
$ws.Range('S15').Value = '25/08/2023 05:12'
$ws.Range('T15').Value = 3.68
$ws.Range('U15').Value = '26/08/2023 16:59'
$ws.Range('V15').Value = 'https://www.betexplorer.com/football/croatia/prva-nl/orijent-vukovar-1991/d0sFnDkj/'

# Row 16
$ws.Range('F16').Value = 'Cibalia'
$ws.Range('G16').Value = 3
$ws.Range('H16').Value = 'Jarun'
$ws.Range('I16').Value = 1
$ws.Range('J16').Value = 1.64
$ws.Range('K16').Value = '25/08/2023 05:12'
$ws.Range('L16').Value = 1.74
$ws.Range('M16').Value = '26/08/2023 16:52'
$ws.Range('N16').Value = 3.71
$ws.Range('O16').Value = '25/08/2023 05:12'
$ws.Range('P16').Value = 3.86
$ws.Range('Q16').Value = '26/08/2023 16:52'
$ws.Range('R16').Value = 4.42
$ws.Range('S16').Value = '25/08/2023 05:12'
$ws.Range('T16').Value = 4.19
$ws.Range('U16').Value = '26/08/2023 16:52'
$ws.Range('V16').Value = 'https://www.betexplorer.com/football/croatia/prva-nl/cibalia-jarun/zitJoX4d/'

# Row 17
$ws.Range('F17').Value = 'Zrinski Jurjevac'
$ws.Range('G17').Value = 3
$ws.Range('H17').Value = 'Solin'
$ws.Range('I17').Value = 0
$ws.Range('J17').Value = 1.95
$ws.Range('K17').Value = '25/08/2023 05:12'
$ws.Range('L17').Value = 2.22
$ws.Range('M17').Value = '26/08/2023 16:56'
$ws.Range('N17').Value = 3.34
$ws.Range('O17').Value = '25/08/2023 05:12'
$ws.Range('P17').Value = 3.32
$ws.Range('Q17').Value = '26/08/2023 16:56'
$ws.Range('R17').Value = 3.29
$ws.Range('S17').Value = '25/08/2023 05:12'
$ws.Range('T17').Value = 3.12
$ws.Range('U17').Value = '26/08/2023 16:56'
$ws.Range('V17').Value = 'https://www.betexplorer.com/football/croatia/prva-nl/zrinski-jurjevac-solin/vPnAmgzp/'

# Row 21
$ws.Range('F21').Value = 'Vukovar 1991'
$ws.Range('G21').Value = 3
$ws.Range('H21').Value = 'Cibalia'
$ws.Range('I21').Value = 1
$ws.Range('J21').Value = 2.07
$ws.Range('K21').Value = '01/09/2023 04:43'
$ws.Range('L21').Value = 1.98
$ws.Range('M21').Value = '02/09/2023 16:28'
$ws.Range('N21').Value = 3.14
$ws.Range('O21').Value = '01/09/2023 04:43'
$ws.Range('P21').Value = 3.32
$ws.Range('Q21').Value = '02/09/2023 16:28'
$ws.Range('R21').Value = 3.3
$ws.Range('S21').Value = '01/09/2023 04:43'
$ws.Range('T21').Value = 3.76
$ws.Range('U21').Value = '02/09/2023 16:28'
$ws.Range('V21').Value = 'https://www.betexplorer.com/football/croatia/prva-nl/vukovar-1991-cibalia/Srlok6Bl/'

# Row 22
$ws.Range('F22').Value = 'Solin'
$ws.Range('G22').Value = 4
$ws.Range('H22').Value = 'Orijent'
$ws.Range('I22').Value = 1
$ws.Range('J22').Value = 2.08
$ws.Range('K22').Value = '01/09/2023 04:43'
$ws.Range('L22').Value = 2.05
$ws.Range('M22').Value = '02/09/2023 16:28'
$ws.Range('N22').Value = 3.34
$ws.Range('O22').Value = '01/09/2023 04:43'
$ws.Range('P22').Value = 3.65
$ws.Range('Q22').Value = '02/09/2023 16:28'
$ws.Range('R22').Value = 2.99
$ws.Range('S22').Value = '01/09/2023 04:43'
$ws.Range('T22').Value = 3.24
$ws.Range('U22').Value = '02/09/2023 16:28'
$ws.Range('V22').Value = 'https://www.betexplorer.com/football/croatia/prva-nl/solin-orijent/A7wjlQQf/'

# Row 24
$ws.Range('F24').Value = 'Bijelo Brdo'
$ws.Range('G24').Value = 0
$ws.Range('H24').Value = 'Croatia Zmijavci'
$ws.Range('I24').Value = 0
$ws.Range('J24').Value = 2.03
$ws.Range('K24').Value = '01/09/2023 04:43'
$ws.Range('L24').Value = 2.42
$ws.Range('M24').Value = '02/09/2023 16:28'
$ws.Range('N24').Value = 3.16
$ws.Range('O24').Value = '01/09/2023 04:43'
$ws.Range('P24').Value = 3.18
$ws.Range('Q24').Value = '02/09/2023 16:28'
$ws.Range('R24').Value = 3.36
$ws.Range('S24').Value = '01/09/2023 04:43'
$ws.Range('T24').Value = 2.91
$ws.Range('U24').Value = '02/09/2023 16:28'
$ws.Range('V24').Value = 'https://www.betexplorer.com/football/croatia/prva-nl/bijelo-brdo-croatia-zmijavci/rgN98fC9/'

# Row 27
$ws.Range('F27').Value = 'Orijent'
$ws.Range('G27').Value = 5
$ws.Range('H27').Value = 'Sibenik'
$ws.Range('I27').Value = 4
$ws.Range('J27').Value = 4.75
$ws.Range('K27').Value = '09/09/2023 14:42'
$ws.Range('L27').Value = 3.89
$ws.Range('M27').Value = '09/09/2023 16:26'
$ws.Range('N27').Value = 3.9
$ws.Range('O27').Value = '09/09/2023 14:42'
$ws.Range('P27').Value = 3.68
$ws.Range('Q27').Value = '09/09/2023 16:23'
$ws.Range('R27').Value = 1.63
$ws.Range('S27').Value = '09/09/2023 14:42'
$ws.Range('T27').Value = 1.84
$ws.Range('U27').Value = '09/09/2023 16:26'
$ws.Range('V27').Value = 'https://www.betexplorer.com/football/croatia/prva-nl/orijent-sibenik/dOu2oOAD/'

# Row 28
$ws.Range('F28').Value = 'Dugopolje'
$ws.Range('G28').Value = 4
$ws.Range('H28').Value = 'Jarun'
$ws.Range('I28').Value = 0
$ws.Range('J28').Value = 1.7
$ws.Range('K28').Value = '08/09/2023 04:42'
$ws.Range('L28').Value = 1.78
$ws.Range('M28').Value = '09/09/2023 16:21'
$ws.Range('N28').Value = 3.62
$ws.Range('O28').Value = '08/09/2023 04:42'
$ws.Range('P28').Value = 3.67
$ws.Range('Q28').Value = '09/09/2023 16:21'
$ws.Range('R28').Value = 3.95
$ws.Range('S28').Value = '08/09/2023 04:42'
$ws.Range('T28').Value = 4.23
$ws.Range('U28').Value = '09/09/2023 16:21'
$ws.Range('V28').Value = 'https://www.betexplorer.com/football/croatia/prva-nl/dugopolje-jarun/Gvap8oIQ/'

# Row 29
$ws.Range('F29').Value = 'Zrinski Jurjevac'
$ws.Range('G29').Value = 1
$ws.Range('H29').Value = 'Bijelo Brdo'
$ws.Range('I29').Value = 0
$ws.Range('J29').Value = 2.1
$ws.Range('K29').Value = '08/09/2023 04:42'
$ws.Range('L29').Value = 1.85
$ws.Range('M29').Value = '09/09/2023 16:23'
$ws.Range('N29').Value = 3.17
$ws.Range('O29').Value = '08/09/2023 04:42'
$ws.Range('P29').Value = 3.4
$ws.Range('Q29').Value = '09/09/2023 16:23'
$ws.Range('R29').Value = 3.2
$ws.Range('S29').Value = '08/09/2023 04:42'
$ws.Range('T29').Value = 4.23
$ws.Range('U29').Value = '09/09/2023 16:23'
$ws.Range('V29').Value = 'https://www.betexplorer.com/football/croatia/prva-nl/zrinski-jurjevac-bijelo-brdo/IXtbn4e7/'

# Row 46
$ws.Range('F46').Value = 'Solin'
$ws.Range('G46').Value = 0
$ws.Range('H46').Value = 'Croatia Zmijavci'
$ws.Range('I46').Value = 0
$ws.Range('J46').Value = 1.89
$ws.Range('K46').Value = '29/09/2023 02:42'
$ws.Range('L46').Value = 1.79
$ws.Range('M46').Value = '30/09/2023 15:29'
$ws.Range('N46').Value = 3.4
$ws.Range('O46').Value = '29/09/2023 02:42'
$ws.Range('P46').Value = 3.71
$ws.Range('Q46').Value = '30/09/2023 15:29'
$ws.Range('R46').Value = 3.41
$ws.Range('S46').Value = '29/09/2023 02:42'
$ws.Range('T46').Value = 4.1
$ws.Range('U46').Value = '30/09/2023 15:29'
$ws.Range('V46').Value = 'https://www.betexplorer.com/football/croatia/prva-nl/solin-croatia-zmijavci/4OeFBKOI/'

# Row 47
$ws.Range('F47').Value = 'Vukovar 1991'
$ws.Range('G47').Value = 2
$ws.Range('H47').Value = 'Dubrava'
$ws.Range('I47').Value = 2
$ws.Range('J47').Value = 1.61
$ws.Range('K47').Value = '29/09/2023 02:42'
$ws.Range('L47').Value = 1.84
$ws.Range('M47').Value = '30/09/2023 15:15'
$ws.Range('N47').Value = 3.7
$ws.Range('O47').Value = '29/09/2023 02:42'
$ws.Range('P47').Value = 3.83
$ws.Range('Q47').Value = '30/09/2023 15:15'
$ws.Range('R47').Value = 4.43
$ws.Range('S47').Value = '29/09/2023 02:42'
$ws.Range('T47').Value = 3.75
$ws.Range('U47').Value = '30/09/2023 15:15'
$ws.Range('V47').Value = 'https://www.betexplorer.com/football/croatia/prva-nl/vukovar-1991-dubrava-zagreb/dbdBC09C/'

# Row 48
$ws.Range('F48').Value = 'Bijelo Brdo'
$ws.Range('G48').Value = 0
$ws.Range('H48').Value = 'Jarun'
$ws.Range('I48').Value = 0
$ws.Range('J48').Value = 1.93
$ws.Range('K48').Value = '29/09/2023 02:42'
$ws.Range('L48').Value = 2.22
$ws.Range('M48').Value = '30/09/2023 15:19'
$ws.Range('N48').Value = 3.38
$ws.Range('O48').Value = '29/09/2023 02:42'
$ws.Range('P48').Value = 3.27
$ws.Range('Q48').Value = '30/09/2023 15:19'
$ws.Range('R48').Value = 3.3
$ws.Range('S48').Value = '29/09/2023 02:42'
$ws.Range('T48').Value = 3.16
$ws.Range('U48').Value = '30/09/2023 15:19'
$ws.Range('V48').Value = 'https://www.betexplorer.com/football/croatia/prva-nl/bijelo-brdo-jarun/I7c7Dtf6/'

# Row 57
$ws.Range('F57').Value = 'Solin'
$ws.Range('G57').Value = 1
$ws.Range('H57').Value = 'Jarun'
$ws.Range('I57').Value = 1
$ws.Range('J57').Value = 1.85
$ws.Range('K57').Value = '13/10/2023 02:13'
$ws.Range('L57').Value = 1.88
$ws.Range('M57').Value = '14/10/2023 14:51'
$ws.Range('N57').Value = 3.54
$ws.Range('O57').Value = '13/10/2023 02:13'
$ws.Range('P57').Value = 3.62
$ws.Range('Q57').Value = '14/10/2023 14:51'
$ws.Range('R57').Value = 3.41
$ws.Range('S57').Value = '13/10/2023 02:13'
$ws.Range('T57').Value = 3.79
$ws.Range('U57').Value = '14/10/2023 14:51'
$ws.Range('V57').Value = 'https://www.betexplorer.com/football/croatia/prva-nl/solin-jarun/6BAb7QlU/'

# Row 58
$ws.Range('F58').Value = 'Zrinski Jurjevac'
$ws.Range('G58').Value = 4
$ws.Range('H58').Value = 'Croatia Zmijavci'
$ws.Range('I58').Value = 0
$ws.Range('J58').Value = 1.68
$ws.Range('K58').Value = '13/10/2023 02:13'
$ws.Range('L58').Value = 1.58
$ws.Range('M58').Value = '14/10/2023 14:53'
$ws.Range('N58').Value = 3.61
$ws.Range('O58').Value = '13/10/2023 02:13'
$ws.Range('P58').Value = 3.86
$ws.Range('Q58').Value = '14/10/2023 14:53'
$ws.Range('R58').Value = 4.26
$ws.Range('S58').Value = '13/10/2023 02:13'
$ws.Range('T58').Value = 5.57
$ws.Range('U58').Value = '14/10/2023 14:53'
$ws.Range('V58').Value = 'https://www.betexplorer.com/football/croatia/prva-nl/zrinski-jurjevac-croatia-zmijavci/dzj8RTZo/'

# Row 59
$ws.Range('F59').Value = 'Cibalia'
$ws.Range('G59').Value = 0
$ws.Range('H59').Value = 'Sesvete'
$ws.Range('I59').Value = 1
$ws.Range('J59').Value = 1.73
$ws.Range('K59').Value = '13/10/2023 02:13'
$ws.Range('L59').Value = 1.75
$ws.Range('M59').Value = '14/10/2023 14:59'
$ws.Range('N59').Value = 3.53
$ws.Range('O59').Value = '13/10/2023 02:13'
$ws.Range('P59').Value = 3.52
$ws.Range('Q59').Value = '14/10/2023 14:59'
$ws.Range('R59').Value = 3.94
$ws.Range('S59').Value = '13/10/2023 02:13'
$ws.Range('T59').Value = 4.64
$ws.Range('U59').Value = '14/10/2023 14:59'
$ws.Range('V59').Value = 'https://www.betexplorer.com/football/croatia/prva-nl/cibalia-sesvete/h2qHP74b/'

# Row 70
$ws.Range('F70').Value = 'Dugopolje'
$ws.Range('G70').Value = 1
$ws.Range('H70').Value = 'Sesvete'
$ws.Range('I70').Value = 1
$ws.Range('J70').Value = 1.72
$ws.Range('K70').Value = '27/10/2023 03:12'
$ws.Range('L70').Value = 1.68
$ws.Range('M70').Value = '27/10/2023 13:14'
$ws.Range('N70').Value = 3.47
$ws.Range('O70').Value = '27/10/2023 03:12'
$ws.Range('P70').Value = 3.76
$ws.Range('Q70').Value = '28/10/2023 14:26'
$ws.Range('R70').Value = 4.07
$ws.Range('S70').Value = '27/10/2023 03:12'
$ws.Range('T70').Value = 4.74
$ws.Range('U70').Value = '28/10/2023 14:26'
$ws.Range('V70').Value = 'https://www.betexplorer.com/football/croatia/prva-nl/dugopolje-sesvete/Yya1NxuS/'

# Row 71
$ws.Range('F71').Value = 'Bijelo Brdo'
$ws.Range('G71').Value = 1
$ws.Range('H71').Value = 'Solin'
$ws.Range('I71').Value = 1
$ws.Range('J71').Value = 2.14
$ws.Range('K71').Value = '27/10/2023 03:12'
$ws.Range('L71').Value = 2.68
$ws.Range('M71').Value = '28/10/2023 14:57'
$ws.Range('N71').Value = 3.17
$ws.Range('O71').Value = '27/10/2023 03:12'
$ws.Range('P71').Value = 2.66
$ws.Range('Q71').Value = '28/10/2023 14:56'
$ws.Range('R71').Value = 3.02
$ws.Range('S71').Value = '27/10/2023 03:12'
$ws.Range('T71').Value = 3.1
$ws.Range('U71').Value = '28/10/2023 14:57'
$ws.Range('V71').Value = 'https://www.betexplorer.com/football/croatia/prva-nl/bijelo-brdo-solin/hM3tS0Qd/'

# --- 2) Append the 4 newly scraped fixtures (rows 88-91) ---

# Row 88
$ws.Range('A87:V87').Copy()
$ws.Range('A88:V88').PasteSpecial(-4122)
$ws.Range('A88').Value = 87
$ws.Range('B88').Value = 'croatia'
$ws.Range('C88').Value = 'prva-nl'
$ws.Range('D88').Value = '2023-2024'
$ws.Range('E88').Value = 45248.5625
$ws.Range('F88').Value = 'Dugopolje'
$ws.Range('G88').Value = 1
$ws.Range('H88').Value = 'Dubrava'
$ws.Range('I88').Value = 2
$ws.Range('J88').Value = 1.81
$ws.Range('K88').Value = '17/11/2023 01:42'
$ws.Range('L88').Value = 2.05
$ws.Range('M88').Value = '18/11/2023 13:21'
$ws.Range('N88').Value = 3.44
$ws.Range('O88').Value = '17/11/2023 01:42'
$ws.Range('P88').Value = 3.39
$ws.Range('Q88').Value = '18/11/2023 13:21'
$ws.Range('R88').Value = 3.81
$ws.Range('S88').Value = '17/11/2023 01:42'
$ws.Range('T88').Value = 3.48
$ws.Range('U88').Value = '18/11/2023 13:21'
$ws.Range('V88').Value = 'https://www.betexplorer.com/football/croatia/prva-nl/dugopolje-dubrava-zagreb/vLCMSXG7/'

# Row 89
$ws.Range('A87:V87').Copy()
$ws.Range('A89:V89').PasteSpecial(-4122)
$ws.Range('A89').Value = 88
$ws.Range('B89').Value = 'croatia'
$ws.Range('C89').Value = 'prva-nl'
$ws.Range('D89').Value = '2023-2024'
$ws.Range('E89').Value = 45248.5625
$ws.Range('F89').Value = 'Orijent'
$ws.Range('G89').Value = 1
$ws.Range('H89').Value = 'Solin'
$ws.Range('I89').Value = 1
$ws.Range('J89').Value = 1.95
$ws.Range('K89').Value = '17/11/2023 01:42'
$ws.Range('L89').Value = 1.88
$ws.Range('M89').Value = '18/11/2023 13:25'
$ws.Range('N89').Value = 3.43
$ws.Range('O89').Value = '17/11/2023 01:42'
$ws.Range('P89').Value = 3.71
$ws.Range('Q89').Value = '18/11/2023 13:25'
$ws.Range('R89').Value = 3.22
$ws.Range('S89').Value = '17/11/2023 01:42'
$ws.Range('T89').Value = 3.67
$ws.Range('U89').Value = '18/11/2023 13:25'
$ws.Range('V89').Value = 'https://www.betexplorer.com/football/croatia/prva-nl/orijent-solin/ljfDkeOK/'

# Row 90
$ws.Range('A87:V87').Copy()
$ws.Range('A90:V90').PasteSpecial(-4122)
$ws.Range('A90').Value = 89
$ws.Range('B90').Value = 'croatia'
$ws.Range('C90').Value = 'prva-nl'
$ws.Range('D90').Value = '2023-2024'
$ws.Range('E90').Value = 45249.5625
$ws.Range('F90').Value = 'Croatia Zmijavci'
$ws.Range('G90').Value = 1
$ws.Range('H90').Value = 'Bijelo Brdo'
$ws.Range('I90').Value = 1
$ws.Range('J90').Value = 1.72
$ws.Range('K90').Value = '18/11/2023 01:42'
$ws.Range('L90').Value = 1.63
$ws.Range('M90').Value = '19/11/2023 13:23'
$ws.Range('N90').Value = 3.47
$ws.Range('O90').Value = '18/11/2023 01:42'
$ws.Range('P90').Value = 3.41
$ws.Range('Q90').Value = '19/11/2023 13:23'
$ws.Range('R90').Value = 4.07
$ws.Range('S90').Value = '18/11/2023 01:42'
$ws.Range('T90').Value = 6.12
$ws.Range('U90').Value = '19/11/2023 13:24'
$ws.Range('V90').Value = 'https://www.betexplorer.com/football/croatia/prva-nl/croatia-zmijavci-bijelo-brdo/0UDITD11/'

# Row 91
$ws.Range('A87:V87').Copy()
$ws.Range('A91:V91').PasteSpecial(-4122)
$ws.Range('A91').Value = 90
$ws.Range('B91').Value = 'croatia'
$ws.Range('C91').Value = 'prva-nl'
$ws.Range('D91').Value = '2023-2024'
$ws.Range('E91').Value = 45249.5625
$ws.Range('F91').Value = 'Zrinski Jurjevac'
$ws.Range('G91').Value = 2
$ws.Range('H91').Value = 'Sibenik'
$ws.Range('I91').Value = 0
$ws.Range('J91').Value = 2.88
$ws.Range('K91').Value = '18/11/2023 01:42'
$ws.Range('L91').Value = 3.22
$ws.Range('M91').Value = '19/11/2023 13:27'
$ws.Range('N91').Value = 3.08
$ws.Range('O91').Value = '18/11/2023 01:42'
$ws.Range('P91').Value = 2.68
$ws.Range('Q91').Value = '19/11/2023 13:27'
$ws.Range('R91').Value = 2.26
$ws.Range('S91').Value = '18/11/2023 01:42'
$ws.Range('T91').Value = 2.58
$ws.Range('U91').Value = '19/11/2023 13:27'
$ws.Range('V91').Value = 'https://www.betexplorer.com/football/croatia/prva-nl/zrinski-jurjevac-sibenik/0QcHlFwR/'

$excel.CutCopyMode = 0

